# Add a new sub-item "5.2.1. Modelo Bigquery" under "5.2. Modelo de Datos"
# on the "Directorio" sheet, growing the Tabla13 table by one data row plus
# the trailing blank row that the source file keeps under the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the formatting of the last existing row (27) onto the two
#        new rows (28 data row, 29 blank spacer row) so no new cell styles
#        are minted in styles.xml.
$ws.Range("A27:H27").Copy()
$ws.Range("A28:A29").PasteSpecial(-4122)   # xlPasteFormats

# Column G on the new data row needs the hyperlink style (same as G26/G25),
# not the plain style copied above.
$ws.Range("G26").Copy()
$ws.Range("G28").PasteSpecial(-4122)       # xlPasteFormats

$excel.CutCopyMode = 0

# --- 2. Fill in the values for the new data row (28), in the same order the
#        original edit introduced the new shared strings (ID, ID_Padre,
#        Nombre, Descripcion, URL) so sharedStrings.xml comes out identical.
$ws.Range("A28").Value = "5.2.1."

# "5.2" must be stored as text (matching the existing "5.2." / "3.1.2" style
# entries), not auto-converted to the number 5.2. Assign it as a text
# formula, then flatten the formula to its cached value via a values-only
# paste so it's stored as a literal shared string with no formula left
# behind.
$ws.Range("B28").Formula = '="5.2"'
$ws.Range("B28").Copy()
$ws.Range("B28").PasteSpecial(-4163)       # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("C28").Value = 3
$ws.Range("D28").Value = "Modelo Bigquery"
$ws.Range("E28").Value = "url"
$ws.Range("F28").Value = "Relacionamiento de dataset con SQL"

# URL + real hyperlink relationship for column G.
$ws.Range("G28").Value = "https://console.cloud.google.com/bigquery?project=modelomarval&ws=!1m4!1m3!3m2!1smodelomarval!2sproyectos"
$ws.Hyperlinks.Add($ws.Range("G28"), "https://console.cloud.google.com/bigquery?project=modelomarval&ws=!1m4!1m3!3m2!1smodelomarval!2sproyectos")

# H28 / row 29 stay blank (already formatted only, from the copy above).

# --- 3. Grow the table (Tabla13) so it covers the new data row and the
#        trailing blank row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H29"))

# --- 4. Match the recorded selection / active cell from the edit.
$null = $ws.Range("D15").Select()
